$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-30 Tuesday", "2024-07-31 Wednesday"),
    @("179÷5=", "778÷6="),
    @("137÷2=", "576÷2="),
    @("856÷9=", "309÷5="),
    @("161÷8=", "825÷8="),
    @("519÷8=", "265÷5="),
    @("106÷2=", "913÷7="),
    @("978÷2=", "300÷5="),
    @("187÷2=", "529÷6="),
    @("487÷4=", "132÷6="),
    @("545÷5=", "573÷5="),
    @("531÷2=", "207÷7="),
    @("558÷9=", "485÷2="),
    @("510÷3=", "402÷2="),
    @("414÷9=", "844÷7="),
    @("439÷4=", "429÷5="),
    @("421÷3=", "786÷6="),
    @("331÷9=", "234÷4="),
    @("154÷9=", "176÷4="),
    @("154÷5=", "389÷8="),
    @("645÷6=", "630÷8="),
    @("495÷8=", "469÷9="),
    @("916÷4=", "411÷3="),
    @("153÷8=", "638÷7="),
    @("461÷9=", "200÷7="),
    @("876÷5=", "830÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
